$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (labels column) and a new row before row 1 (header row)
$ws.Columns.Item(1).Insert()
$ws.Rows.Item(1).Insert()

# Set the new column A width (target stored width 54.450195 characters)
$ws.Columns.Item(1).ColumnWidth = 53.666667

# Header row
$ws.Range("B1").Value = "Valid"
$ws.Range("C1").Value = "T"
$ws.Range("D1").Value = "Z"
$ws.Range("E1").Value = "p-value"

# Row labels (column A) for each data row
$ws.Range("A2").Value = "CyclomaticComplexity(CC) & CyclomaticComplexity(CC)"
$ws.Range("A3").Value = "CyclomaticComplexity(CC) & NbOperators"
$ws.Range("A4").Value = "MaintainabilityIndex & MaintainabilityIndex"
$ws.Range("A5").Value = "NbOperands & NbOperands"
$ws.Range("A6").Value = "NbUniqueOperators & NbUniqueOperators"
$ws.Range("A7").Value = "NbOperators & CyclomaticComplexity(CC)"
$ws.Range("A8").Value = "NbOperators & NbOperators"
$ws.Range("A9").Value = "ProgramLength & ProgramLength"
$ws.Range("A10").Value = "VocabularySize & VocabularySize"
$ws.Range("A11").Value = "ProgramVolume & ProgramVolume"
$ws.Range("A12").Value = "ProgramLevel & ProgramLevel"
$ws.Range("A13").Value = "EffortToImplement & EffortToImplement"
$ws.Range("A14").Value = "TimeToImplement & TimeToImplement"

